# VT_SHR_ELC_V01.xlsx edit script
# - Update the absolute-path-derived workbook bookmark is not reachable via
#   the Excel object model, so it is left untouched by this script.
# - Change several numeric values on the ProcessCharac sheet so they differ
#   from what the students already have, clear a handful of cells so the
#   students have to work the right values out themselves, and mark the
#   edited/cleared cells with fill colors (green = changed, yellow = cleared)
#   so they stand out.
# - Move the active sheet/tab selection from SEC_Comm to ProcessCharac
#   (cell P13 selected there) to mirror where the instructor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProcessCharac")

# ---- Row 12 : FT-ELCCOA / ETCOASTM1E --------------------------------------
$ws.Range("H12").Value = 0.421
$ws.Range("M12").Value = 0.872
$ws.Range("N12").Value = 261.4

# ---- Row 13 : ECCOABP1E ----------------------------------------------------
$ws.Range("H13").Value = 0.34
$ws.Range("M13").Value = 0.85
$ws.Range("N13").ClearContents()
$ws.Range("P13").ClearContents()

# ---- Row 15 : ETNGAGT1E ----------------------------------------------------
$ws.Range("H15").Value = 0.44
$ws.Range("M15").Value = 0.91
$ws.Range("N15").Value = 237

# ---- Row 16 : ERWINON1E ----------------------------------------------------
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 562

# ---- Highlight cells that were merely cleared (students fill these in) ----
$ws.Range("P13").Interior.Color = 65535
$ws.Range("M16").Interior.Color = 65535

# ---- Highlight cells whose numbers were changed (efficiency column) -------
$ws.Range("H12").Interior.Color = 5296274
$ws.Range("H13").Interior.Color = 5296274
$ws.Range("H15").Interior.Color = 5296274

# ---- Highlight cells whose numbers were changed (AFA / capacity columns) --
$ws.Range("M12").Interior.Color = 5296274
$ws.Range("N12").Interior.Color = 5296274
$ws.Range("M13").Interior.Color = 5296274
$ws.Range("M15").Interior.Color = 5296274
$ws.Range("N15").Interior.Color = 5296274
$ws.Range("N16").Interior.Color = 5296274

# ---- N13 was cleared too, but keep it formatted as an integer -------------
$ws.Range("N13").NumberFormat = "0"
$ws.Range("N13").Interior.Color = 65535

# ---- Move the selection / active tab to ProcessCharac!P13 -----------------
$ws.Range("P13").Select()
$ws.Activate()
